$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Reward -> salaries, 9000 -> 100000, date -> 45782.125185185185
$ws.Range("A2").Value = "salaries"
$ws.Range("B2").Value = 100000
$ws.Range("C2").Value = 45782.125185185185

# Row 3: Salary -> "freelance " (trailing space), 20000 -> 50000, date -> 45782.125185185185
$ws.Range("A3").Value = "freelance "
$ws.Range("B3").Value = 50000
$ws.Range("C3").Value = 45782.125185185185

# Row 4: Freelance -> crypto, 6000 -> 200000, date -> 45782.125185185185
$ws.Range("A4").Value = "crypto"
$ws.Range("B4").Value = 200000
$ws.Range("C4").Value = 45782.125185185185
